# Version 2-8-17 Finalizado filtro Clases de Material
# Insert a new column before column B on the "Gasto Capital" sheet, label the
# new header cell "CODIGO" (merged across the two header rows like the other
# header cells), and give it the same look (font/border/alignment) as the
# neighbouring header cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gasto Capital")

# Shift everything from column B onward one column to the right.
$ws.Columns.Item(2).Insert()

# New header cell for the inserted column.
$ws.Range("B9").Value = "CODIGO"
$ws.Range("B9:B10").Merge()

# Match the formatting of the adjoining header cell (font, border, alignment).
$ws.Range("C9:C10").Copy()
$ws.Range("B9:B10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Leave the cursor where the author's session ended up.
$ws.Range("C13").Select() | Out-Null
